$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update summary statistics for the "CrudeOil" column after the
# underlying data manipulation (min value changed from -37.63 to 10.01,
# which shifts the mean/std for CrudeOil and the rounding-level std
# values for the other columns).

$ws.Range("B3").Value = 63.51691476969101   # mean (CrudeOil)
$ws.Range("B4").Value = 25.8285177638053    # std  (CrudeOil)
$ws.Range("C4").Value = 509.6281056181531   # std  (Gold)
$ws.Range("D4").Value = 96.25514278043239   # std  (S&P500)
$ws.Range("E4").Value = 1037.069513617562   # std  (FTSE)
$ws.Range("B5").Value = 10.01000022888184   # min  (CrudeOil)
